# Updated TPM values for Lgi2-Adam23 LR-pairs sheet (new sending-cluster rotation:
# FAPs / MuSCs / Resolving-Mac instead of ECs / FAPs / MuSCs).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = New-Object "object[,]" 1,20
$row2[0,0] = "FAPs"
$row2[0,1] = "Lgi2"
$row2[0,2] = "Adam23"
$row2[0,3] = "ECs"
$row2[0,4] = 3
$row2[0,5] = 1
$row2[0,6] = 10.27464766666667
$row2[0,7] = 30.823943
$row2[0,8] = 0.9550701914680073
$row2[0,9] = 0.9550701914680074
$row2[0,10] = 3
$row2[0,11] = 1
$row2[0,12] = 0.165747
$row2[0,13] = 0.497241
$row2[0,14] = 0.008095785894995438
$row2[0,15] = 0.00809578589499544
$row2[0,16] = 1.702992026807
$row2[0,17] = 15.326928241263
$row2[0,18] = 0.007732043784817286
$row2[0,19] = 0.007732043784817288
$ws.Range("A2:T2").Value = $row2

$row3 = New-Object "object[,]" 1,20
$row3[0,0] = "FAPs"
$row3[0,1] = "Lgi2"
$row3[0,2] = "Adam23"
$row3[0,3] = "FAPs"
$row3[0,4] = 3
$row3[0,5] = 1
$row3[0,6] = 10.27464766666667
$row3[0,7] = 30.823943
$row3[0,8] = 0.9550701914680073
$row3[0,9] = 0.9550701914680074
$row3[0,10] = 3
$row3[0,11] = 1
$row3[0,12] = 14.494489
$row3[0,13] = 43.483467
$row3[0,14] = 0.7079722685862583
$row3[0,15] = 0.7079722685862583
$row3[0,16] = 148.9257675833757
$row3[0,17] = 1340.331908250381
$row3[0,18] = 0.6761632101127172
$row3[0,19] = 0.6761632101127173
$ws.Range("A3:T3").Value = $row3

$row4 = New-Object "object[,]" 1,20
$row4[0,0] = "FAPs"
$row4[0,1] = "Lgi2"
$row4[0,2] = "Adam23"
$row4[0,3] = "MuSCs"
$row4[0,4] = 3
$row4[0,5] = 1
$row4[0,6] = 10.27464766666667
$row4[0,7] = 30.823943
$row4[0,8] = 0.9550701914680073
$row4[0,9] = 0.9550701914680074
$row4[0,10] = 3
$row4[0,11] = 1
$row4[0,12] = 5.642879333333333
$row4[0,13] = 16.928638
$row4[0,14] = 0.2756221404547972
$row4[0,15] = 0.2756221404547972
$row4[0,16] = 57.97859697551489
$row4[0,17] = 521.807372779634
$row4[0,18] = 0.2632384904569851
$row4[0,19] = 0.2632384904569852
$ws.Range("A4:T4").Value = $row4

$row5 = New-Object "object[,]" 1,20
$row5[0,0] = "FAPs"
$row5[0,1] = "Lgi2"
$row5[0,2] = "Adam23"
$row5[0,3] = "Resolving-Mac"
$row5[0,4] = 3
$row5[0,5] = 1
$row5[0,6] = 10.27464766666667
$row5[0,7] = 30.823943
$row5[0,8] = 0.9550701914680073
$row5[0,9] = 0.9550701914680074
$row5[0,10] = 3
$row5[0,11] = 1
$row5[0,12] = 0.1701286666666667
$row5[0,13] = 0.510386
$row5[0,14] = 0.008309805063949155
$row5[0,15] = 0.008309805063949155
$row5[0,16] = 1.748012107999778
$row5[0,17] = 15.732108971998
$row5[0,18] = 0.007936447113487736
$row5[0,19] = 0.007936447113487738
$ws.Range("A5:T5").Value = $row5

$row6 = New-Object "object[,]" 1,20
$row6[0,0] = "MuSCs"
$row6[0,1] = "Lgi2"
$row6[0,2] = "Adam23"
$row6[0,3] = "ECs"
$row6[0,4] = 3
$row6[0,5] = 1
$row6[0,6] = 0.475652
$row6[0,7] = 1.426956
$row6[0,8] = 0.04421378342596928
$row6[0,9] = 0.04421378342596929
$row6[0,10] = 3
$row6[0,11] = 1
$row6[0,12] = 0.165747
$row6[0,13] = 0.497241
$row6[0,14] = 0.008095785894995438
$row6[0,15] = 0.00809578589499544
$row6[0,16] = 0.07883789204399999
$row6[0,17] = 0.709541028396
$row6[0,18] = 0.0003579453242243452
$row6[0,19] = 0.0003579453242243453
$ws.Range("A6:T6").Value = $row6

$row7 = New-Object "object[,]" 1,20
$row7[0,0] = "MuSCs"
$row7[0,1] = "Lgi2"
$row7[0,2] = "Adam23"
$row7[0,3] = "FAPs"
$row7[0,4] = 3
$row7[0,5] = 1
$row7[0,6] = 0.475652
$row7[0,7] = 1.426956
$row7[0,8] = 0.04421378342596928
$row7[0,9] = 0.04421378342596929
$row7[0,10] = 3
$row7[0,11] = 1
$row7[0,12] = 14.494489
$row7[0,13] = 43.483467
$row7[0,14] = 0.7079722685862583
$row7[0,15] = 0.7079722685862583
$row7[0,16] = 6.894332681828001
$row7[0,17] = 62.048994136452
$row7[0,18] = 0.03130213255486498
$row7[0,19] = 0.03130213255486498
$ws.Range("A7:T7").Value = $row7

$row8 = New-Object "object[,]" 1,20
$row8[0,0] = "MuSCs"
$row8[0,1] = "Lgi2"
$row8[0,2] = "Adam23"
$row8[0,3] = "MuSCs"
$row8[0,4] = 3
$row8[0,5] = 1
$row8[0,6] = 0.475652
$row8[0,7] = 1.426956
$row8[0,8] = 0.04421378342596928
$row8[0,9] = 0.04421378342596929
$row8[0,10] = 3
$row8[0,11] = 1
$row8[0,12] = 5.642879333333333
$row8[0,13] = 16.928638
$row8[0,14] = 0.2756221404547972
$row8[0,15] = 0.2756221404547972
$row8[0,16] = 2.684046840658666
$row8[0,17] = 24.156421565928
$row8[0,18] = 0.01218629762547049
$row8[0,19] = 0.01218629762547049
$ws.Range("A8:T8").Value = $row8

$row9 = New-Object "object[,]" 1,20
$row9[0,0] = "MuSCs"
$row9[0,1] = "Lgi2"
$row9[0,2] = "Adam23"
$row9[0,3] = "Resolving-Mac"
$row9[0,4] = 3
$row9[0,5] = 1
$row9[0,6] = 0.475652
$row9[0,7] = 1.426956
$row9[0,8] = 0.04421378342596928
$row9[0,9] = 0.04421378342596929
$row9[0,10] = 3
$row9[0,11] = 1
$row9[0,12] = 0.1701286666666667
$row9[0,13] = 0.510386
$row9[0,14] = 0.008309805063949155
$row9[0,15] = 0.008309805063949155
$row9[0,16] = 0.08092204055733333
$row9[0,17] = 0.728298365016
$row9[0,18] = 0.0003674079214094707
$row9[0,19] = 0.0003674079214094708
$ws.Range("A9:T9").Value = $row9

$row10 = New-Object "object[,]" 1,20
$row10[0,0] = "Resolving-Mac"
$row10[0,1] = "Lgi2"
$row10[0,2] = "Adam23"
$row10[0,3] = "ECs"
$row10[0,4] = 1
$row10[0,5] = 0.3333333333333333
$row10[0,6] = 0.007703000000000001
$row10[0,7] = 0.023109
$row10[0,8] = 0.0007160251060233982
$row10[0,9] = 0.0007160251060233983
$row10[0,10] = 3
$row10[0,11] = 1
$row10[0,12] = 0.165747
$row10[0,13] = 0.497241
$row10[0,14] = 0.008095785894995438
$row10[0,15] = 0.00809578589499544
$row10[0,16] = 0.001276749141
$row10[0,17] = 0.011490742269
$row10[0,18] = 0.000005796785953806840703485172
$row10[0,19] = 0.000005796785953806842397551067
$ws.Range("A10:T10").Value = $row10

$row11 = New-Object "object[,]" 1,20
$row11[0,0] = "Resolving-Mac"
$row11[0,1] = "Lgi2"
$row11[0,2] = "Adam23"
$row11[0,3] = "FAPs"
$row11[0,4] = 1
$row11[0,5] = 0.3333333333333333
$row11[0,6] = 0.007703000000000001
$row11[0,7] = 0.023109
$row11[0,8] = 0.0007160251060233982
$row11[0,9] = 0.0007160251060233983
$row11[0,10] = 3
$row11[0,11] = 1
$row11[0,12] = 14.494489
$row11[0,13] = 43.483467
$row11[0,14] = 0.7079722685862583
$row11[0,15] = 0.7079722685862583
$row11[0,16] = 0.111651048767
$row11[0,17] = 1.004859438903
$row11[0,18] = 0.0005069259186761013
$row11[0,19] = 0.0005069259186761014
$ws.Range("A11:T11").Value = $row11

$row12 = New-Object "object[,]" 1,20
$row12[0,0] = "Resolving-Mac"
$row12[0,1] = "Lgi2"
$row12[0,2] = "Adam23"
$row12[0,3] = "MuSCs"
$row12[0,4] = 1
$row12[0,5] = 0.3333333333333333
$row12[0,6] = 0.007703000000000001
$row12[0,7] = 0.023109
$row12[0,8] = 0.0007160251060233982
$row12[0,9] = 0.0007160251060233983
$row12[0,10] = 3
$row12[0,11] = 1
$row12[0,12] = 5.642879333333333
$row12[0,13] = 16.928638
$row12[0,14] = 0.2756221404547972
$row12[0,15] = 0.2756221404547972
$row12[0,16] = 0.04346709950466667
$row12[0,17] = 0.391203895542
$row12[0,18] = 0.0001973523723415421
$row12[0,19] = 0.0001973523723415421
$ws.Range("A12:T12").Value = $row12

$row13 = New-Object "object[,]" 1,20
$row13[0,0] = "Resolving-Mac"
$row13[0,1] = "Lgi2"
$row13[0,2] = "Adam23"
$row13[0,3] = "Resolving-Mac"
$row13[0,4] = 1
$row13[0,5] = 0.3333333333333333
$row13[0,6] = 0.007703000000000001
$row13[0,7] = 0.023109
$row13[0,8] = 0.0007160251060233982
$row13[0,9] = 0.0007160251060233983
$row13[0,10] = 3
$row13[0,11] = 1
$row13[0,12] = 0.1701286666666667
$row13[0,13] = 0.510386
$row13[0,14] = 0.008309805063949155
$row13[0,15] = 0.008309805063949155
$row13[0,16] = 0.001310501119333334
$row13[0,17] = 0.011794510074
$row13[0,18] = 0.00000595002905194796478602088
$row13[0,19] = 0.000005950029051947965633053827
$ws.Range("A13:T13").Value = $row13

